$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-11 (columns A, B, D, E). Column C stays blank.
$data = @(
    @{Row=2;  A='S253441'; B=1; D=2; E='E2004703EC9060269CAC0110'},
    @{Row=3;  A='S253441'; B=1; D=3; E='E2004702ED6060268CB9010D'},
    @{Row=4;  A='S234145'; B=1; D=8; E='E20047053EC06026B1CF0108'},
    @{Row=5;  A='S454132'; B=7; D=3; E='E2004704D9C06026AB7F0114'},
    @{Row=6;  A='S534241'; B=5; D=5; E='E20047053EA06026B1CD010A'},
    @{Row=7;  A='S243415'; B=2; D=3; E='E2004704D9906026AB7C010D'},
    @{Row=8;  A='S342451'; B=5; D=8; E='E2004704D9B06026AB7E0109'},
    @{Row=9;  A='S542314'; B=3; D=1; E='E2004704D9A06026AB7D010E'},
    @{Row=10; A='S452341'; B=2; D=1; E='E20047053EB06026B1CE010A'},
    @{Row=11; A='S412354'; B=4; D=7; E='E20047053E906026B1CC0105'}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
